$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the class-name (B column) labels to the new set.
$ws.Range("B2").Value = "管实2001"
$ws.Range("B3").Value = "管工2001"
$ws.Range("B4").Value = "管工2002"
$ws.Range("B5").Value = "管工2003"
$ws.Range("B6").Value = "管工2004"
$ws.Range("B7").Value = "管实1901"
$ws.Range("B8").Value = "信管1901"
$ws.Range("B9").Value = "信管1902"
$ws.Range("B10").Value = "物流1901"
$ws.Range("B11").Value = "物流1902"

# Move the active selection from D10 to D11.
$ws.Range("D11").Select()
